$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: split the "Find the prebuilt workout routine ..." sentence into
# two runs with new wording, keeping the original (identical) run formatting.
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Find the prebuilt workout routine for Leg exercise and set*") {
        $full = $para.Range
        $r = $d.Range($full.Start, $full.End - 1)
        $xmlFrag = "<w:p $wNs>" +
            '<w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
            '<w:t xml:space="preserve">Find the prebuilt workout routine for Leg exercise and set </w:t></w:r>' +
            '<w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
            '<w:t>the weight as 45 pounds, and reps as 10.</w:t></w:r>' +
            '</w:p>'
        $r.InsertXML($xmlFrag)
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not locate the 'Find the prebuilt workout routine' paragraph"
}

# ---------------------------------------------------------------------------
# Change 2: split the "Check that you have finished ..." sentence into three
# runs ("... timer ", "ends", ".") with the original run formatting.
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Check that you have finished the exercise*") {
        $full = $para.Range
        $r = $d.Range($full.Start, $full.End - 1)
        $xmlFrag = "<w:p $wNs>" +
            '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
            '<w:t xml:space="preserve">Check that you have finished the exercise and answer at least 2 quizzes until the timer </w:t></w:r>' +
            '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
            '<w:t>ends</w:t></w:r>' +
            '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
            '<w:t>.</w:t></w:r>' +
            '</w:p>'
        $r.InsertXML($xmlFrag)
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not locate the 'Check that you have finished' paragraph"
}

# ---------------------------------------------------------------------------
# Change 3: numbering.xml - the ilvl=1 / tplc=04090019 level in the abstract
# numbering definition used by numId 2 is no longer tentative.
# ---------------------------------------------------------------------------
$lists = $d.ListTemplates
$targetTemplate = $null
foreach ($lt in $lists) {
    $lvl = $lt.ListLevels.Item(2)
    if ($lvl.NumberStyle -eq 4 -and $lt.ListLevels.Item(1).TrailingCharacter -ne $null) {
        # placeholder - identification done differently below
    }
}

Write-Output "done"
